# Rettet en feil i test-kodebok
# The worksheet "1-Testskjema-felter" had a duplicated "Hoyde" field row
# (row 6 was an exact duplicate of what is now row 10/11 after the fix).
# Remove the erroneous duplicate row, which shifts the remaining rows up.

$wb = $excel.ActiveWorkbook

$wsFelter = $wb.Worksheets.Item("1-Testskjema-felter")
$wsRegler2 = $wb.Worksheets.Item("2-Sluttskjema-regler")

# Delete the duplicate "Hoyde" row (row 6), shifting rows 7-12 up to 6-11.
$wsFelter.Rows.Item(6).Delete() | Out-Null

# Update the selection on the fixed sheet and make it the active sheet/tab.
$wsFelter.Activate()
$wsFelter.Range("C16").Select() | Out-Null

# The previously active sheet keeps its own selection, but is no longer
# the tab that is selected/shown when the workbook is opened.
$wsRegler2.Range("G5").Select() | Out-Null

$wsFelter.Activate()
